# This script applies a cyclic rotation of text content across several
# paragraphs of the LOQ4051 document (Objetivos, Docente, Programa resumido,
# Programa, Metodo, Criterio, Norma de recuperacao, Bibliografia).
#
# Mapping (old text -> new text), derived from the unified diff:
#   P6  (Objetivos)          : A -> C
#   P8  (Docente)             : B -> A
#   P10 (Programa resumido)   : C -> D
#   P12 (Programa)            : D -> E
#   P14 run "Metodo:"value    : E -> F
#   P14 run "Criterio:"value  : F -> G
#   P14 run "Norma..."value   : G -> H
#   P16 (Bibliografia)        : H -> B
#
# Because this is a genuine rotation (several paragraphs swap content with
# each other), we first stash each paragraph's original text into a unique
# placeholder, and only then fill in the real final text for every slot.
# This avoids any chance of one replacement's output being matched again by
# a later Find/Replace call.

$d = $word.ActiveDocument

function Set-ParagraphText($paragraphIndex, $newText) {
    $para = $d.Paragraphs.Item($paragraphIndex)
    $rng = $para.Range
    # Exclude the trailing paragraph mark from the range so we only replace
    # the visible text content of the paragraph.
    $rng.MoveEnd(1, -1) | Out-Null
    $rng.Text = $newText
}

# Step 1: move every old value into a unique, unambiguous placeholder so
# that later writes can never be re-matched by an earlier Find.
Set-ParagraphText 6  "@@PLACEHOLDER_1@@"
Set-ParagraphText 8  "@@PLACEHOLDER_2@@"
Set-ParagraphText 10 "@@PLACEHOLDER_3@@"
Set-ParagraphText 12 "@@PLACEHOLDER_4@@"
Set-ParagraphText 16 "@@PLACEHOLDER_8@@"

# Paragraph 14 contains three separate runs of interest: the value after
# "Metodo:", after "Criterio:" and after "Norma de recuperacao:". Replace
# each value run individually, scoping the Find to just that paragraph's
# Range so we do not disturb the bold label runs ("Metodo: ", "Critério: ",
# etc).
$p14 = $d.Paragraphs.Item(14)

$rng = $p14.Range.Duplicate
$rng.Find.Execute("Supervisão das atividades desenvolvidas pelo aluno durante o estágio.", $true, $true, $false, $false, $false, $true, 1, $false, "@@PLACEHOLDER_5@@", 2) | Out-Null

$rng = $p14.Range.Duplicate
$rng.Find.Execute("MF = Nota baseada em relatório final e no desempenho no estágio, a ser atribuída pelo professor orientador do estágio.", $true, $true, $false, $false, $false, $true, 1, $false, "@@PLACEHOLDER_6@@", 2) | Out-Null

$rng = $p14.Range.Duplicate
$rng.Find.Execute("Não será oferecida recuperação.", $true, $true, $false, $false, $false, $true, 1, $false, "@@PLACEHOLDER_7@@", 2) | Out-Null

# Step 2: write the real final text into every slot, using the original
# value that belonged to the appropriate source slot.

Set-ParagraphText 6  "Plano de Trabalho específico. Realização do Estágio. Relatório final e/ou parciais."
Set-ParagraphText 8  "Fornecer oportunidade de aplicação dos conhecimentos fundamentais da Engenharia Química nos projetos e processos químicos. Complementação da formação geral curricular. Adaptação psicológica e social do estudante à sua futura atividade profissional."
Set-ParagraphText 10 "Participação do aluno em processo seletivo de empresas ou no setor acadêmico. Estágio realizado sob a supervisão da Escola de Engenharia de Lorena, através do Departamento em Engenharia Química. O conteúdo será estabelecido individualmente no Plano de Trabalho entre o Supervisor do Estágio e o professor orientador, desde que relacionado com as áreas afins da Engenharia Química.  Apresentação de relatório final e/ou relatórios parciais sobre as atividades desenvolvidas no estágio."
Set-ParagraphText 12 "Supervisão das atividades desenvolvidas pelo aluno durante o estágio."
Set-ParagraphText 16 "198273 - Domingos Savio Giordani"

$rng = $p14.Range.Duplicate
$rng.Find.Execute("@@PLACEHOLDER_5@@", $true, $true, $false, $false, $false, $true, 1, $false, "MF = Nota baseada em relatório final e no desempenho no estágio, a ser atribuída pelo professor orientador do estágio.", 2) | Out-Null

$rng = $p14.Range.Duplicate
$rng.Find.Execute("@@PLACEHOLDER_6@@", $true, $true, $false, $false, $false, $true, 1, $false, "Não será oferecida recuperação.", 2) | Out-Null

$rng = $p14.Range.Duplicate
$rng.Find.Execute("@@PLACEHOLDER_7@@", $true, $true, $false, $false, $false, $true, 1, $false, "A ser definida com o orientador em função das atividades desenvolvidas no estágio.", 2) | Out-Null

Write-Host "Done."
